$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 41: lec80 follow-up question (written before the C35 edit so the
# shared-string table ends up in the same append order as the target file)
$ws.Range("C40").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Value = "lec80, co jest bardziej czasochlonne for do polowy rozmiaru tablicy ale dodatkowa linijka kodu czy krotszy kod ale for po calej tablicy"
$ws.Rows.Item(41).RowHeight = 28.5

# Extend the question text in C35 and grow its row height to fit
$ws.Range("C35").Value = "Można się odwoływać w metodzie Main do innych metod z innych klas instancja.metoda. Jak się można odwołać, nie będąc w metodzie Main do innej klasy? - ogólnie jak najlepiej się odwolywac do  metod z innych klas czy pakietów."
$ws.Rows.Item(35).RowHeight = 42.75

# New row 42: scanner.nextLine question
$ws.Range("C40").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C42").Value = "scanner.nextLine() - ""to clear the input buffer' ? Po co jest ta komenda?"

# New row 43: date + lec82 Main note
$ws.Range("B35").Copy()
$ws.Range("B43").PasteSpecial(-4122)
$ws.Range("B43").Value = 43295

$ws.Range("C40").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("C43").Value = "lec82 Main. "

$ws.Application.CutCopyMode = $false

# View state updates: scroll so row 27 is at the top and select the new last cell
$ws.Application.Goto($ws.Range("A27"))
$ws.Range("C43").Select()
